# Add two new rows (116 and 117) of data to each of the 6 worksheets.
# New dates: 2025-11-24 (serial 45985) and 2025-11-25 (serial 45986)

$wb = $excel.ActiveWorkbook

# Values per sheet (sheet index -> [[date, amount], [date, amount]])
$data = @{
    1 = @(@(45985, 518214), @(45986, 514718))
    2 = @(@(45985, 350958), @(45986, 329310))
    3 = @(@(45985, 151151), @(45986, 147616))
    4 = @(@(45985, 214010), @(45986, 217840))
    5 = @(@(45985, 689087), @(45986, 709767))
    6 = @(@(45985, 75099), @(45986, 83334))
}

foreach ($sheetIndex in 1..6) {
    $ws = $wb.Worksheets.Item($sheetIndex)

    $rows = $data[$sheetIndex]

    $startRow = 116
    for ($i = 0; $i -lt $rows.Count; $i++) {
        $rowNum = $startRow + $i
        $dateVal = $rows[$i][0]
        $amtVal = $rows[$i][1]

        $cellA = $ws.Cells.Item($rowNum, 1)
        $cellB = $ws.Cells.Item($rowNum, 2)

        # Copy the number format / style from the row above so the new
        # date cell matches the existing date column formatting.
        $prevA = $ws.Cells.Item($rowNum - 1, 1)
        $cellA.NumberFormat = $prevA.NumberFormat

        $cellA.Value = $dateVal
        $cellB.Value = $amtVal
    }
}
